# Add new columns I (I0) and J (IF) to Sheet1, mirroring the style of
# existing header cells and filling in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy style from H1 so the new headers match existing ones
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-13 for columns I and J (both columns share the same value per row)
$values = @{
    2  = 6
    3  = 6
    4  = 7
    5  = 8
    6  = 8
    7  = 9
    8  = 8
    9  = 5
    10 = 8
    11 = 8
    12 = 7
    13 = 8
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v   # column I
    $ws.Cells.Item($row, 10).Value = $v  # column J
}
